$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the Ciboulette data block (rows 552-553),
# pushing the existing rows 552:627 down to 554:629.
$ws.Rows("552:553").Insert()

# New row 552: Primera quality, week of 2023-02-27 (serial 44984)
$ws.Range("A552").Value2 = 6
$ws.Range("B552").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C552").Value2 = "Metropolitana"
$ws.Range("D552").Value2 = 44984
$ws.Range("E552").Value2 = 13
$ws.Range("F552").Value2 = 100112039
$ws.Range("G552").Value2 = "Ciboulette"
$ws.Range("H552").Value2 = "Sin especificar"
$ws.Range("I552").Value2 = "Primera"
$ws.Range("J552").Value2 = 250
$ws.Range("K552").Value2 = 2000
$ws.Range("L552").Value2 = 2000
$ws.Range("M552").Value2 = 2000
$ws.Range("N552").Value2 = "`$/docena de atados"
$ws.Range("O552").Value2 = "Región Metropolitana"
$ws.Range("P552").Value2 = 667
$ws.Range("Q552").Value2 = 3
$ws.Range("R552").Value2 = "Hortaliza"

# New row 553: Segunda quality, same week (serial 44984)
$ws.Range("A553").Value2 = 6
$ws.Range("B553").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C553").Value2 = "Metropolitana"
$ws.Range("D553").Value2 = 44984
$ws.Range("E553").Value2 = 13
$ws.Range("F553").Value2 = 100112039
$ws.Range("G553").Value2 = "Ciboulette"
$ws.Range("H553").Value2 = "Sin especificar"
$ws.Range("I553").Value2 = "Segunda"
$ws.Range("J553").Value2 = 120
$ws.Range("K553").Value2 = 1500
$ws.Range("L553").Value2 = 1500
$ws.Range("M553").Value2 = 1500
$ws.Range("N553").Value2 = "`$/docena de atados"
$ws.Range("O553").Value2 = "Región Metropolitana"
$ws.Range("P553").Value2 = 500
$ws.Range("Q553").Value2 = 3
$ws.Range("R553").Value2 = "Hortaliza"
